$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '66.091.86'
Set-TextValue "E2" '  +0.81%  '
Set-TextValue "D3" '3.184.27'
Set-TextValue "E3" '  -4.44%  '
Set-TextValue "E4" '  -0.02%  '
Set-TextValue "D5" '572.68'
Set-TextValue "E5" '  -0.53%  '
Set-TextValue "D6" '172.87'
Set-TextValue "E6" '  -3.49%  '
Set-TextValue "B7" 'XRP'
Set-TextValue "C7" 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue "D7" '0.605'
Set-TextValue "E7" '  -2.48%  '
Set-TextValue "B8" 'USDC'
Set-TextValue "C8" 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue "D8" '1.00'
Set-TextValue "E8" '  -0.04%  '
Set-TextValue "D9" '3.178.34'
Set-TextValue "E9" '  -4.55%  '
Set-TextValue "D10" '0.127'
Set-TextValue "E10" '  -0.73%  '
Set-TextValue "D11" '6.63'
Set-TextValue "E11" '  -3.32%  '
Set-TextValue "D12" '0.396'
Set-TextValue "E12" '  -2.53%  '
Set-TextValue "D13" '3.734.90'
Set-TextValue "E13" '  -4.45%  '
Set-TextValue "D14" '0.136'
Set-TextValue "E14" '  +1.04%  '
Set-TextValue "D15" '27.31'
Set-TextValue "E15" '  -4.05%  '
Set-TextValue "D16" '66.002.18'
Set-TextValue "E16" '  +0.76%  '
Set-TextValue "D17" '0.0000165'
Set-TextValue "E17" '  -1.86%  '
Set-TextValue "D18" '3.179.38'
Set-TextValue "E18" '  -5.53%  '
Set-TextValue "D19" '5.76'
Set-TextValue "E19" '  -0.16%  '
Set-TextValue "D20" '13.01'
Set-TextValue "E20" '  -3.01%  '
Set-TextValue "D21" '364.77'
Set-TextValue "E21" '  +0.22%  '
Set-TextValue "D22" '7.31'
Set-TextValue "E22" '  -1.60%  '
Set-TextValue "D23" '1.00'
Set-TextValue "E23" '  +0.15%  '
Set-TextValue "D24" '69.21'
Set-TextValue "E24" '  -3.19%  '
Set-TextValue "D25" '0.499'
Set-TextValue "E25" '  -3.74%  '
Set-TextValue "D26" '3.318.86'
Set-TextValue "E26" '  -4.79%  '
Set-TextValue "E27" '  -5.73%  '
Set-TextValue "D28" '9.93'
Set-TextValue "E28" '  +2.87%  '
Set-TextValue "D29" '0.177'
Set-TextValue "E29" '  -0.18%  '
Set-TextValue "E30" '  +0.04%  '
Set-TextValue "B31" 'USDe'
Set-TextValue "C31" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D31" '0.998'
Set-TextValue "E31" '  -0.12%  '
Set-TextValue "B32" 'PancakeSwap'
Set-TextValue "C32" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D32" '1.93'
Set-TextValue "E32" '  -1.27%  '
Set-TextValue "D33" '5.42'
Set-TextValue "E33" '  -4.37%  '
Set-TextValue "D34" '22.14'
Set-TextValue "E34" '  -3.11%  '
Set-TextValue "D35" '6.63'
Set-TextValue "E35" '  -2.56%  '
Set-TextValue "E36" '  -1.66%  '
Set-TextValue "D37" '163.43'
Set-TextValue "E37" '  +2.22%  '
Set-TextValue "E38" '  -0.93%  '
Set-TextValue "D39" '0.829'
Set-TextValue "E39" '  -2.00%  '
Set-TextValue "D40" '1.81'
Set-TextValue "E40" '  +3.63%  '
Set-TextValue "D41" '26.26'
Set-TextValue "E41" '  -3.42%  '
Set-TextValue "D42" '2.54'
Set-TextValue "E42" '  +0.29%  '
Set-TextValue "D43" '2.658.70'
Set-TextValue "E43" '  -2.52%  '
Set-TextValue "D44" '6.21'
Set-TextValue "E44" '  -0.30%  '
Set-TextValue "E45" '  -1.75%  '
Set-TextValue "D46" '40.01'
Set-TextValue "E46" '  +0.18%  '
Set-TextValue "E47" '  -0.87%  '
Set-TextValue "D48" '328.18'
Set-TextValue "E48" '  -1.45%  '
Set-TextValue "D49" '23.95'
Set-TextValue "E49" '  -0.38%  '
Set-TextValue "D50" '0.0275'
Set-TextValue "E50" '  -0.76%  '
Set-TextValue "E51" '  -1.79%  '